# Fixed issues #5 #6 #7 #10 #12
# Relocate the "marker" labels (n <LABEL> n) within the game board grid on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: move the "SHP" marker from D3:F3 to J3:L3, relabeled as "HSE"
$ws.Range("D3").Value = "e"
$ws.Range("E3").Value = "e"
$ws.Range("F3").Value = "e"
$ws.Range("J3").Value = "n"
$ws.Range("K3").Value = "HSE"
$ws.Range("L3").Value = "n"

# Row 5: move the "FAC" marker from D5:F5 to J5:L5, relabeled as "HWY"
$ws.Range("D5").Value = "e"
$ws.Range("E5").Value = "e"
$ws.Range("F5").Value = "e"
$ws.Range("J5").Value = "n"
$ws.Range("K5").Value = "HWY"
$ws.Range("L5").Value = "n"

# Row 5: add a new "SHP" marker at P5:R5
$ws.Range("P5").Value = "n"
$ws.Range("Q5").Value = "SHP"
$ws.Range("R5").Value = "n"

# Row 5: add a new "HSE" marker at V5:X5
$ws.Range("V5").Value = "n"
$ws.Range("W5").Value = "HSE"
$ws.Range("X5").Value = "n"

# Row 7: add a new "BCH" marker at J7:L7
$ws.Range("J7").Value = "n"
$ws.Range("K7").Value = "BCH"
$ws.Range("L7").Value = "n"

# Row 7: add a new "FAC" marker at V7:X7
$ws.Range("V7").Value = "n"
$ws.Range("W7").Value = "FAC"
$ws.Range("X7").Value = "n"

# Row 9: add a new "BCH" marker at V9:X9
$ws.Range("V9").Value = "n"
$ws.Range("W9").Value = "BCH"
$ws.Range("X9").Value = "n"
